$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 2950.8333
$ws.Cells.Item(98, 9).Value = 3041
$ws.Cells.Item(98, 10).Value = 2500
$ws.Cells.Item(98, 11).Value = 3041
$ws.Cells.Item(98, 12).Value = 2500
$ws.Cells.Item(98, 13).Value = -1543
$ws.Cells.Item(98, 14).Value = -5496
$ws.Cells.Item(106, 8).Value = 6098.0625
$ws.Cells.Item(106, 9).Value = 5837.933
$ws.Cells.Item(106, 11).Value = 5837.933
$ws.Cells.Item(106, 13).Value = -5206.933
$ws.Cells.Item(107, 8).Value = 1452.1052
$ws.Cells.Item(107, 9).Value = 1860.8462
$ws.Cells.Item(107, 11).Value = 1860.8462
$ws.Cells.Item(107, 13).Value = 59.15380000000005
$ws.Cells.Item(122, 8).Value = 2950.8333
$ws.Cells.Item(122, 9).Value = 3041
$ws.Cells.Item(122, 10).Value = 2500
$ws.Cells.Item(122, 11).Value = 9123
$ws.Cells.Item(122, 12).Value = 7500
$ws.Cells.Item(122, 13).Value = -6673
$ws.Cells.Item(122, 14).Value = -12400
$ws.Cells.Item(137, 8).Value = 2251.024
$ws.Cells.Item(137, 9).Value = 2024.9706
$ws.Cells.Item(137, 11).Value = 6074.9118
$ws.Cells.Item(137, 13).Value = -3524.9118

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4852.9873
$ws.Cells.Item(32, 9).Value = 4967.338
$ws.Cells.Item(32, 11).Value = 4967.338
$ws.Cells.Item(32, 13).Value = -4680.338
$ws.Cells.Item(61, 8).Value = 2720.875
$ws.Cells.Item(61, 9).Value = 2545.3333
$ws.Cells.Item(61, 11).Value = 2545.3333
$ws.Cells.Item(61, 13).Value = -2333.3333
$ws.Cells.Item(122, 8).Value = 4485.722
$ws.Cells.Item(122, 9).Value = 3961.05
$ws.Cells.Item(122, 10).Value = 5141.5625
$ws.Cells.Item(122, 11).Value = 11883.15
$ws.Cells.Item(122, 12).Value = 15424.6875
$ws.Cells.Item(122, 13).Value = -9433.150000000001
$ws.Cells.Item(122, 14).Value = -20324.6875
$ws.Cells.Item(132, 8).Value = 1939.8518
$ws.Cells.Item(132, 9).Value = 2125.0527
$ws.Cells.Item(132, 11).Value = 6375.158100000001
$ws.Cells.Item(132, 13).Value = -3845.158100000001
$ws.Cells.Item(136, 8).Value = 2720.875
$ws.Cells.Item(136, 9).Value = 2545.3333
$ws.Cells.Item(136, 11).Value = 7635.999899999999
$ws.Cells.Item(136, 13).Value = -5085.999899999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 348.33334
$ws.Cells.Item(80, 10).Value = 339.2
$ws.Cells.Item(80, 12).Value = 339.2
$ws.Cells.Item(80, 14).Value = -2335.2
$ws.Cells.Item(83, 8).Value = 348.33334
$ws.Cells.Item(83, 10).Value = 339.2
$ws.Cells.Item(83, 12).Value = 1696
$ws.Cells.Item(83, 14).Value = -11680
$ws.Cells.Item(107, 8).Value = 3117.1904
$ws.Cells.Item(107, 9).Value = 2818.9143
$ws.Cells.Item(107, 11).Value = 2818.9143
$ws.Cells.Item(107, 13).Value = -898.9142999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1503.7097
$ws.Cells.Item(31, 9).Value = 1393.1482
$ws.Cells.Item(31, 11).Value = 1393.1482
$ws.Cells.Item(31, 13).Value = -1098.1482
$ws.Cells.Item(34, 8).Value = 1503.7097
$ws.Cells.Item(34, 9).Value = 1393.1482
$ws.Cells.Item(34, 11).Value = 1393.1482
$ws.Cells.Item(34, 13).Value = -1191.1482
$ws.Cells.Item(58, 8).Value = 1472.4286
$ws.Cells.Item(58, 9).Value = 852.63635
$ws.Cells.Item(58, 10).Value = 2521.3076
$ws.Cells.Item(58, 11).Value = 852.63635
$ws.Cells.Item(58, 12).Value = 2521.3076
$ws.Cells.Item(58, 13).Value = -649.63635
$ws.Cells.Item(58, 14).Value = -2927.3076
$ws.Cells.Item(86, 8).Value = 32275.6
$ws.Cells.Item(86, 9).Value = 42819.875
$ws.Cells.Item(86, 10).Value = 20225
$ws.Cells.Item(86, 11).Value = 42819.875
$ws.Cells.Item(86, 12).Value = 20225
$ws.Cells.Item(86, 13).Value = -41696.875
$ws.Cells.Item(86, 14).Value = -22471
$ws.Cells.Item(89, 8).Value = 32275.6
$ws.Cells.Item(89, 9).Value = 42819.875
$ws.Cells.Item(89, 10).Value = 20225
$ws.Cells.Item(89, 11).Value = 214099.375
$ws.Cells.Item(89, 12).Value = 101125
$ws.Cells.Item(89, 13).Value = -208483.375
$ws.Cells.Item(89, 14).Value = -112357
$ws.Cells.Item(121, 8).Value = 15000
$ws.Cells.Item(121, 9).Value = 15000
$ws.Cells.Item(121, 10).Value = 0
$ws.Cells.Item(121, 11).Value = 15000
$ws.Cells.Item(121, 12).Value = 0
$ws.Cells.Item(121, 13).Value = -13690
$ws.Cells.Item(136, 8).Value = 1472.4286
$ws.Cells.Item(136, 9).Value = 852.63635
$ws.Cells.Item(136, 10).Value = 2521.3076
$ws.Cells.Item(136, 11).Value = 2557.90905
$ws.Cells.Item(136, 12).Value = 7563.9228
$ws.Cells.Item(136, 13).Value = -7.909050000000207
$ws.Cells.Item(136, 14).Value = -12663.9228
$ws.Cells.Item(121, 14).ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 48362060
$ws.Cells.Item(4, 9).Value = 56422370
$ws.Cells.Item(4, 11).Value = 169267110
$ws.Cells.Item(4, 13).Value = -169266998
$ws.Cells.Item(5, 8).Value = 1034.6875
$ws.Cells.Item(5, 9).Value = 1020.6071
$ws.Cells.Item(5, 10).Value = 1133.25
$ws.Cells.Item(5, 11).Value = 3061.8213
$ws.Cells.Item(5, 12).Value = 3399.75
$ws.Cells.Item(5, 13).Value = -2949.8213
$ws.Cells.Item(5, 14).Value = -3623.75
$ws.Cells.Item(7, 8).Value = 5535.409
$ws.Cells.Item(7, 9).Value = 8943.154
$ws.Cells.Item(7, 10).Value = 613.1111
$ws.Cells.Item(7, 11).Value = 26829.462
$ws.Cells.Item(7, 12).Value = 1839.3333
$ws.Cells.Item(7, 13).Value = -26717.462
$ws.Cells.Item(7, 14).Value = -2063.3333
$ws.Cells.Item(13, 8).Value = 191.47058
$ws.Cells.Item(13, 9).Value = 227.57143
$ws.Cells.Item(13, 11).Value = 682.71429
$ws.Cells.Item(13, 13).Value = -514.71429
$ws.Cells.Item(107, 8).Value = 2099.724
$ws.Cells.Item(107, 10).Value = 1970.1
$ws.Cells.Item(107, 12).Value = 5910.299999999999
$ws.Cells.Item(107, 14).Value = -9750.299999999999
$ws.Cells.Item(121, 8).Value = 1460.9166
$ws.Cells.Item(121, 10).Value = 1790.2858
$ws.Cells.Item(121, 12).Value = 5370.857400000001
$ws.Cells.Item(121, 14).Value = -7990.857400000001
$ws.Cells.Item(135, 8).Value = 1034.6875
$ws.Cells.Item(135, 9).Value = 1020.6071
$ws.Cells.Item(135, 10).Value = 1133.25
$ws.Cells.Item(135, 11).Value = 9185.463899999999
$ws.Cells.Item(135, 12).Value = 10199.25
$ws.Cells.Item(135, 13).Value = -6650.463899999999
$ws.Cells.Item(135, 14).Value = -15269.25
$ws.Cells.Item(138, 8).Value = 6346.147
$ws.Cells.Item(138, 9).Value = 4603
$ws.Cells.Item(138, 10).Value = 7895.6113
$ws.Cells.Item(138, 11).Value = 13809
$ws.Cells.Item(138, 12).Value = 23686.8339
$ws.Cells.Item(138, 13).Value = -8669
$ws.Cells.Item(138, 14).Value = -33966.8339

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 2833.1667
$ws.Cells.Item(46, 9).Value = 1360
$ws.Cells.Item(46, 10).Value = 3885.4285
$ws.Cells.Item(46, 11).Value = 1360
$ws.Cells.Item(46, 12).Value = 3885.4285
$ws.Cells.Item(46, 13).Value = -1172
$ws.Cells.Item(46, 14).Value = -4261.4285
$ws.Cells.Item(47, 8).Value = 37247.25
$ws.Cells.Item(47, 9).Value = 0
$ws.Cells.Item(47, 10).Value = 37247.25
$ws.Cells.Item(47, 11).Value = 0
$ws.Cells.Item(47, 12).Value = 37247.25
$ws.Cells.Item(47, 14).Value = -38227.25
$ws.Cells.Item(52, 8).Value = 37247.25
$ws.Cells.Item(52, 9).Value = 0
$ws.Cells.Item(52, 10).Value = 37247.25
$ws.Cells.Item(52, 11).Value = 0
$ws.Cells.Item(52, 12).Value = 37247.25
$ws.Cells.Item(52, 14).Value = -37713.25
$ws.Cells.Item(132, 8).Value = 2086.9644
$ws.Cells.Item(132, 9).Value = 1793.8889
$ws.Cells.Item(132, 11).Value = 5381.6667
$ws.Cells.Item(132, 13).Value = -2851.6667
$ws.Cells.Item(136, 8).Value = 2225.439
$ws.Cells.Item(136, 9).Value = 2086.1516
$ws.Cells.Item(136, 10).Value = 2800
$ws.Cells.Item(136, 11).Value = 6258.4548
$ws.Cells.Item(136, 12).Value = 8400
$ws.Cells.Item(136, 13).Value = -3708.4548
$ws.Cells.Item(136, 14).Value = -13500
$ws.Cells.Item(47, 13).ClearContents()
$ws.Cells.Item(52, 13).ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(75, 8).Value = 75000
$ws.Cells.Item(75, 10).Value = 75000
$ws.Cells.Item(75, 12).Value = 75000
$ws.Cells.Item(75, 14).Value = -76872
$ws.Cells.Item(78, 8).Value = 75000
$ws.Cells.Item(78, 10).Value = 75000
$ws.Cells.Item(78, 12).Value = 225000
$ws.Cells.Item(78, 14).Value = -234360
$ws.Cells.Item(132, 8).Value = 2187.3845
$ws.Cells.Item(132, 9).Value = 1344.1
$ws.Cells.Item(132, 11).Value = 4032.3
$ws.Cells.Item(132, 13).Value = -1502.3
